# Generate Report for Handoff
# Adds a new tracked file "e20b7ece-5702-47f7-8036-18b450605ec8" as row 9
# on the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$fileId = "e20b7ece-5702-47f7-8036-18b450605ec8"
$mdName = "$fileId.md"
$zhXlf  = "$fileId.48e0524c51e165e7493e667c5bae51a5cd410726.zh-cn.xlf"
$deXlf  = "$fileId.48e0524c51e165e7493e667c5bae51a5cd410726.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = $mdName
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-32-14 06:32:09"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e20b7ece5702eeee57025702eeee570257025702/e2e/$mdName",
    "",
    "",
    $mdName
)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A9").Value = $mdName
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = $zhXlf
$wsZhCn.Range("E9").Value = "2016-03-14 06:32:07"
$wsZhCn.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I9").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e20b7ece5702eeee57025702eeee570257025702/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e20b7ece5702eeee57025702eeee570257025702/e2e/$mdName",
    "",
    "",
    ".md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e20b7ece5702eeee57025702eeee570257025702/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf",
    "",
    "",
    $zhXlf
)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A9").Value = $mdName
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = $deXlf
$wsDeDe.Range("E9").Value = "2016-03-14 06:32:09"
$wsDeDe.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I9").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e20b7ece5702eeee57025702eeee570257025702/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e20b7ece5702eeee57025702eeee570257025702/e2e/$mdName",
    "",
    "",
    ".md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e20b7ece5702eeee57025702eeee570257025702/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf",
    "",
    "",
    $deXlf
)
